# Applies the "Natmi following Dr Hou advice" update: adds the ECs sending
# cluster (full 3x3 Sending x Target grid) and refreshes the recomputed
# edge-weight statistics for rows 2-10 of the LR-pairs sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @("ECs", "ECs", "ECs", "FAPs", "FAPs", "FAPs", "sCs", "sCs", "sCs")
$colB = @("Angpt1", "Angpt1", "Angpt1", "Angpt1", "Angpt1", "Angpt1", "Angpt1", "Angpt1", "Angpt1")
$colC = @("Itgb1", "Itgb1", "Itgb1", "Itgb1", "Itgb1", "Itgb1", "Itgb1", "Itgb1", "Itgb1")
$colD = @("ECs", "FAPs", "sCs", "ECs", "FAPs", "sCs", "ECs", "FAPs", "sCs")
$colE = @(1, 1, 1, 3, 3, 3, 3, 3, 3)
$colF = @(0.3333333333333333, 0.3333333333333333, 0.3333333333333333, 1, 1, 1, 1, 1, 1)
$colG = @(0.07171233333333334, 0.07171233333333334, 0.07171233333333334, 11.05178533333333, 11.05178533333333, 11.05178533333333, 3.370524666666666, 3.370524666666666, 3.370524666666666)
$colH = @(0.215137, 0.215137, 0.215137, 33.155356, 33.155356, 33.155356, 10.111574, 10.111574, 10.111574)
$colI = @(0.004947717871829782, 0.004947717871829782, 0.004947717871829782, 0.7625064374239614, 0.7625064374239614, 0.7625064374239614, 0.2325458447042087, 0.2325458447042087, 0.2325458447042087)
$colJ = @(0.004947717871829783, 0.004947717871829783, 0.004947717871829783, 0.7625064374239615, 0.7625064374239615, 0.7625064374239615, 0.2325458447042088, 0.2325458447042088, 0.2325458447042088)
$colK = @(3, 3, 3, 3, 3, 3, 3, 3, 3)
$colL = @(1, 1, 1, 1, 1, 1, 1, 1, 1)
$colM = @(112.513392, 106.314466, 124.6916553333333, 112.513392, 106.314466, 124.6916553333333, 112.513392, 106.314466, 124.6916553333333)
$colN = @(337.540176, 318.943398, 374.074966, 337.540176, 318.943398, 374.074966, 337.540176, 318.943398, 374.074966)
$colO = @(0.3275312977368564, 0.3094859589441663, 0.3629827433189773, 0.3275312977368564, 0.3094859589441663, 0.3629827433189773, 0.3275312977368564, 0.3094859589441663, 0.3629827433189773)
$colP = @(0.3275312977368564, 0.3094859589441664, 0.3629827433189773, 0.3275312977368564, 0.3094859589441664, 0.3629827433189773, 0.3275312977368564, 0.3094859589441664, 0.3629827433189773)
$colQ = @(8.068597871568, 7.624058423947333, 8.941929551149112, 1243.473855509184, 1174.964656059965, 1378.065407601988, 379.229163066336, 358.3355300764946, 420.2763000284982)
$colR = @(72.61738084411199, 68.616525815526, 80.477365960342, 11191.26469958265, 10574.68190453969, 12402.58866841789, 3413.062467597023, 3225.019770688452, 3782.486700256484)
$colS = @(0.001620532455396246, 0.00153124921014843, 0.001795936206285106, 0.2497447229821771, 0.2359850359872546, 0.2767766784545296, 0.07616604229928295, 0.07196967374676322, 0.08441012865816255)
$colT = @(0.001620532455396246, 0.00153124921014843, 0.001795936206285107, 0.2497447229821772, 0.2359850359872547, 0.2767766784545297, 0.07616604229928296, 0.07196967374676325, 0.08441012865816258)

for ($i = 0; $i -lt $colA.Count; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 1).Value = $colA[$i]
  $ws.Cells.Item($r, 2).Value = $colB[$i]
  $ws.Cells.Item($r, 3).Value = $colC[$i]
  $ws.Cells.Item($r, 4).Value = $colD[$i]
  $ws.Cells.Item($r, 5).Value = $colE[$i]
  $ws.Cells.Item($r, 6).Value = $colF[$i]
  $ws.Cells.Item($r, 7).Value = $colG[$i]
  $ws.Cells.Item($r, 8).Value = $colH[$i]
  $ws.Cells.Item($r, 9).Value = $colI[$i]
  $ws.Cells.Item($r, 10).Value = $colJ[$i]
  $ws.Cells.Item($r, 11).Value = $colK[$i]
  $ws.Cells.Item($r, 12).Value = $colL[$i]
  $ws.Cells.Item($r, 13).Value = $colM[$i]
  $ws.Cells.Item($r, 14).Value = $colN[$i]
  $ws.Cells.Item($r, 15).Value = $colO[$i]
  $ws.Cells.Item($r, 16).Value = $colP[$i]
  $ws.Cells.Item($r, 17).Value = $colQ[$i]
  $ws.Cells.Item($r, 18).Value = $colR[$i]
  $ws.Cells.Item($r, 19).Value = $colS[$i]
  $ws.Cells.Item($r, 20).Value = $colT[$i]
}
